$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.515.11'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '1.626.91'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '213.12'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").Value = '0.504'
$ws.Range("E6").Value = '  +1.91%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.249'
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.0623'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '18.80'
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("D12").Value = '1.853.25'
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("D13").Value = '1.634.51'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("D16").Value = '65.12'
$ws.Range("E16").Value = '  +3.28%  '
$ws.Range("D17").Value = '26.537.38'
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").Value = '214.63'
$ws.Range("E19").Value = '  +2.63%  '
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("D22").Value = '6.26'
$ws.Range("E22").Value = '  +1.34%  '
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("D24").Value = '2.10'
$ws.Range("E24").Value = '  +10.09%  '
$ws.Range("D25").Value = '147.69'
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("E28").Value = '  +2.03%  '
$ws.Range("D29").Value = '15.53'
$ws.Range("E29").Value = '  +1.05%  '
$ws.Range("E30").Value = '  -1.66%  '
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = '1.243.00'
$ws.Range("E34").Value = '  +6.32%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("E37").Value = '  +4.38%  '
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("D39").Value = '0.509'
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("D41").Value = '2.27'
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("E42").Value = '  +0.58%  '
$ws.Range("D43").Value = '5.33'
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("D44").Value = '1.763.36'
$ws.Range("E44").Value = '  -0.85%  '
$ws.Range("D45").Value = '93.25'
$ws.Range("E45").Value = '  +1.26%  '
$ws.Range("E46").Value = '  +2.01%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0104'
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '54.86'
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.0958'
$ws.Range("E50").Value = '  +2.10%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.407'
$ws.Range("E51").Value = '  -0.63%  '
